$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
# for each data row (2-23). Derived from the canonical diff.
$data = @{
    2  = @(44175, 1600, 1000, 1200, 1100, 1100)
    3  = @(44883, 800,  550,  600,  575,  575)
    4  = @(44607, 900,  1300, 1400, 1350, 1350)
    5  = @(44449, 1300, 900,  950,  925,  925)
    6  = @(44687, 1000, 1200, 1300, 1250, 1250)
    7  = @(44407, 1000, 1200, 1300, 1250, 1250)
    8  = @(44341, 1300, 900,  1000, 950,  950)
    9  = @(44453, 1000, 800,  900,  850,  850)
    10 = @(44284, 1500, 800,  850,  825,  825)
    11 = @(44229, 1500, 1400, 1500, 1450, 1450)
    12 = @(44673, 900,  1300, 1400, 1350, 1350)
    13 = @(44442, 1250, 850,  900,  875,  875)
    14 = @(44784, 1000, 1200, 1300, 1250, 1250)
    15 = @(44476, 900,  700,  800,  750,  750)
    16 = @(44638, 1000, 900,  950,  925,  925)
    17 = @(44649, 600,  900,  1000, 950,  950)
    18 = @(44243, 1200, 1200, 1300, 1250, 1250)
    19 = @(44550, 1300, 1000, 1200, 1100, 1100)
    20 = @(44656, 1000, 900,  1000, 950,  950)
    21 = @(44455, 1100, 900,  1000, 950,  950)
    22 = @(44484, 900,  750,  800,  775,  775)
    23 = @(44291, 1000, 1000, 1200, 1100, 1100)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]   # D - Fecha
    $ws.Cells.Item($row, 10).Value = $vals[1]  # J - Volumen
    $ws.Cells.Item($row, 11).Value = $vals[2]  # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals[3]  # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals[4]  # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals[5]  # P - Precio $/Kg
}
